$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("matched")

$ws2 = $wb.Worksheets.Item("unmatched_invoices")
$ws2.Cells.Clear()
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("A1").Value = "info"
$ws2.Range("A2").Value = "No data"

$ws3 = $wb.Worksheets.Item("unmatched_payments")
$ws3.Cells.Clear()
$ws1.Range("A1").Copy()
$ws3.Range("A1").PasteSpecial(-4122)
$ws3.Range("A1").Value = "info"
$ws3.Range("A2").Value = "No data"
